# "Added more test cases"
# Fill in the "passed" status for the existing test-case rows 11-15, then
# append three brand-new test cases (rows 16-18) with their statuses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows 11-15 get a "passed" status in column B.
$ws.Range("B11").Value = "passed"
$ws.Range("B12").Value = "passed"
$ws.Range("B13").Value = "passed"
$ws.Range("B14").Value = "passed"
$ws.Range("B15").Value = "passed"

# New test case rows.
$ws.Range("A16").Value = "Searched apps show up in search box"
$ws.Range("B16").Value = "passed"

$ws.Range("A17").Value = "new users receive confirmation email"
$ws.Range("B17").Value = "passed"

$ws.Range("A18").Value = "Admins can view a page to view newly submitted apps"

# Match the updated view state: scrolled down a bit, with A19 selected as
# the next empty row to fill in.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
[void]$ws.Range("A19").Select()
